{"js": "// \"added login page2 (window1)\" \u2014 append short \" done\" / \" \" notes to five\n// existing lines (Prathmesh / Aniket / Saurab(h) / prathmesh / saurabh),\n// each as a brand-new trailing run so the original run is left untouched.\n\nconst body = context.document.body;\n\nasync function appendRun(searchText, newText) {\n  const results = body.search(searchText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + searchText);\n  }\n\n  // Use the first match; every search string below is unique in this document.\n  const found = results.items[0];\n\n  // Collapse to the end of the found range, then insert new text there.\n  // Inserting via a collapsed \"End\" range (instead of calling insertText\n  // directly on the match, or using getRange(\"After\")) keeps the newly\n  // typed text in its own run instead of merging it back into the run\n  // that was just searched.\n  const endRange = found.getRange(\"End\");\n  endRange.insertText(newText, \"End\");\n  await context.sync();\n}\n\n// 1) \"... AUTHORITY  Prathmesh\" -> \"... AUTHORITY  Prathmesh done\"\nawait appendRun(\" Prathmesh\", \" done\");\n\n// 2) \"USERS  Aniket\" -> \"USERS  Aniket \" (trailing space run added)\nawait appendRun(\" Aniket\", \" \");\n\n// 3) \"RESOURCES Saurab\" -> \"RESOURCES Saurabh done\"\nawait appendRun(\"Saurab\", \"h done\");\n\n// 4) \"Window 1 -prathmesh\" -> \"Window 1 -prathmesh done\"\nawait appendRun(\"prathmesh\", \" done\");\n\n// 5) \"Winow 3  saurabh\" -> \"Winow 3  saurabh done\"\nawait appendRun(\"saurabh\", \" done\");\n", "ps1": "# \"added login page2 (window1)\" \u2014 append short \" done\" / \" \" notes to five\n# existing lines (Prathmesh / Aniket / Saurab(h) / prathmesh / saurabh),\n# each as a brand-new trailing run so the original run is left untouched.\n\n$d = $word.ActiveDocument\n\nfunction Append-AfterText($findText, $newText) {\n    # Search the whole document body each time (fresh Range) so a prior\n    # insertion can't leave a stale/collapsed Range behind.\n    $rng = $d.Content\n    $found = $rng.Find.Execute($findText, $true)   # FindText, MatchCase:=True\n    if (-not $found) {\n        throw \"Could not find text: $findText\"\n    }\n    # $rng now spans exactly the matched text; collapse to its end point\n    # and insert there so the new text becomes its own run, immediately\n    # after (not merged into) the run that was matched.\n    $rng.Collapse($wdCollapseEnd)\n    $rng.InsertAfter($newText)\n}\n\n# 1) \"... AUTHORITY  Prathmesh\" -> \"... AUTHORITY  Prathmesh done\"\nAppend-AfterText \" Prathmesh\" \" done\"\n\n# 2) \"USERS  Aniket\" -> \"USERS  Aniket \" (trailing space run added)\nAppend-AfterText \" Aniket\" \" \"\n\n# 3) \"RESOURCES Saurab\" -> \"RESOURCES Saurabh done\"\nAppend-AfterText \"Saurab\" \"h done\"\n\n# 4) \"Window 1 -prathmesh\" -> \"Window 1 -prathmesh done\"\nAppend-AfterText \"prathmesh\" \" done\"\n\n# 5) \"Winow 3  saurabh\" -> \"Winow 3  saurabh done\"\nAppend-AfterText \"saurabh\" \" done\"\n"}
